# Generate Report for Handoff
# Rename the source markdown file's generated GUID from
# ce65fb08-6f80-43f5-bfeb-48667ac90013 to 95a67d60-01bc-4fc6-9c32-acf04af32e97,
# and refresh the zh-cn / de-de handoff package names + handoff timestamps
# produced by the new run.

$wb = $excel.ActiveWorkbook

$oldGuidMd   = "ce65fb08-6f80-43f5-bfeb-48667ac90013.md"
$newGuidMd   = "95a67d60-01bc-4fc6-9c32-acf04af32e97.md"

$oldZhCn     = "ce65fb08-6f80-43f5-bfeb-48667ac90013.59d1009a497c9f3e1e812eb0b772449af7836785.zh-cn.xlf"
$newZhCn     = "95a67d60-01bc-4fc6-9c32-acf04af32e97.40060195a60408424decca15bff4093c42d0980d.zh-cn.xlf"

$oldDeDe     = "ce65fb08-6f80-43f5-bfeb-48667ac90013.59d1009a497c9f3e1e812eb0b772449af7836785.de-de.xlf"
$newDeDe     = "95a67d60-01bc-4fc6-9c32-acf04af32e97.40060195a60408424decca15bff4093c42d0980d.de-de.xlf"

$oldZhCnTime = "2016-02-24 09:50:51"
$newZhCnTime = "2016-02-24 09:51:50"

$oldDeDeTime = "2016-02-24 09:51:06"
$newDeDeTime = "2016-02-24 09:52:01"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # ---- update cell values (shared strings) ----
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $v = $cell.Value()
            if ($v -eq $oldGuidMd) { $cell.Value = $newGuidMd }
            elseif ($v -eq $oldZhCn) { $cell.Value = $newZhCn }
            elseif ($v -eq $oldDeDe) { $cell.Value = $newDeDe }
            elseif ($v -eq $oldZhCnTime) { $cell.Value = $newZhCnTime }
            elseif ($v -eq $oldDeDeTime) { $cell.Value = $newDeDeTime }
        }
    }

    # ---- update hyperlink display text to match (leave Address/r:id alone) ----
    foreach ($h in $ws.Hyperlinks) {
        $disp = $h.TextToDisplay
        if ($disp -eq $oldGuidMd) { $h.TextToDisplay = $newGuidMd }
        elseif ($disp -eq $oldZhCn) { $h.TextToDisplay = $newZhCn }
        elseif ($disp -eq $oldDeDe) { $h.TextToDisplay = $newDeDe }
    }
}
